# Generate Report for handback
# Applies:
#  1. Status text "Not yet handed off" -> "Handed back" everywhere it appears
#     (Overview!B2, Overview!C2, zh-cn!B2, de-de!B2)
#  2. zh-cn sheet (row 2): add "Latest Target File" (E2) and
#     "Latest Handback File" (F2) hyperlinked cells, and set the
#     "Latest Handback DateTime" (G2) to the handback timestamp.
#  3. de-de sheet (row 2): set "Latest Handoff Datetime" (D2), add
#     "Latest Target File" (E2) / "Latest Handback File" (F2) hyperlinked
#     cells, and set "Latest Handback DateTime" (G2).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back"

# 1. Status flips from "Not yet handed off" to "Handed back" for every
#    sheet / column that surfaces it.
$wsOverview.Range("B2").Value = $handedBack
$wsOverview.Range("C2").Value = $handedBack
$wsZhCn.Range("B2").Value = $handedBack
$wsDeDe.Range("B2").Value = $handedBack

# Hyperlink targets reused from the already-present links in each sheet.
$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/2c60ad50c60e89aed7e0240939d70735a63928a6/e2e/b3bea00a-d71f-490e-a515-aee0e91a7f98.md"
$zhCnXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4b3778732febbcbbcaf324272f317a1b29dc490b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/b3bea00a-d71f-490e-a515-aee0e91a7f98.3d9bee042f38885cc9d3a056ecd147b35e867d8b.zh-cn.xlf"
$deDeXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f566126f10a52ad0931bff2092bad5f572f2b217/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/b3bea00a-d71f-490e-a515-aee0e91a7f98.3d9bee042f38885cc9d3a056ecd147b35e867d8b.de-de.xlf"

$mdFileName = "b3bea00a-d71f-490e-a515-aee0e91a7f98.md"
$zhCnXlfFileName = "b3bea00a-d71f-490e-a515-aee0e91a7f98.3d9bee042f38885cc9d3a056ecd147b35e867d8b.zh-cn.xlf"
$deDeXlfFileName = "b3bea00a-d71f-490e-a515-aee0e91a7f98.3d9bee042f38885cc9d3a056ecd147b35e867d8b.de-de.xlf"

$hyperlinkColor = 15570276

# 2. zh-cn (handback completed at 2016-01-08 15:08:44)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), $mdTarget, "", "", $mdFileName)
$wsZhCn.Range("E2").Font.Color = $hyperlinkColor

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhCnXlfTarget, "", "", $zhCnXlfFileName)
$wsZhCn.Range("F2").Font.Color = $hyperlinkColor

$wsZhCn.Range("G2").Value = "2016-01-08 15:08:44"

# 3. de-de (handback completed at 2016-01-08 15:09:07)
$wsDeDe.Range("D2").Value = "2016-01-08 15:06:49"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), $mdTarget, "", "", $mdFileName)
$wsDeDe.Range("E2").Font.Color = $hyperlinkColor

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deDeXlfTarget, "", "", $deDeXlfFileName)
$wsDeDe.Range("F2").Font.Color = $hyperlinkColor

$wsDeDe.Range("G2").Value = "2016-01-08 15:09:07"
